$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 65; this shifts the existing rows
# 65..150 down to 66..151 (matching the dimension change A1:R150 -> A1:R151).
$ws.Rows("65:65").Insert()

# Populate the newly inserted (blank) row 65 with the new record's data.
$ws.Range("A65").Value2 = 8
$ws.Range("B65").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C65").Value2 = "Coquimbo"
$ws.Range("D65").Value2 = 44413
$ws.Range("E65").Value2 = 4
$ws.Range("F65").Value2 = 100112032
$ws.Range("G65").Value2 = "Zapallo italiano"
$ws.Range("H65").Value2 = "Sin especificar"
$ws.Range("I65").Value2 = "Primera"
$ws.Range("J65").Value2 = 500
$ws.Range("K65").Value2 = 7500
$ws.Range("L65").Value2 = 8000
$ws.Range("M65").Value2 = 7750
$ws.Range("N65").Value2 = "`$/caja 50 unidades"
$ws.Range("O65").Value2 = "Región de Arica y Parinacota"
$ws.Range("P65").Value2 = 155
$ws.Range("Q65").Value2 = 50
$ws.Range("R65").Value2 = "Hortaliza"
